$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1557.375
$ws.Range("I6").Value = 1260
$ws.Range("J6").Value = 2053
$ws.Range("K6").Value = 3780
$ws.Range("L6").Value = 6159
$ws.Range("M6").Value = -3668
$ws.Range("N6").Value = -6383

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 649.4
$ws.Range("I28").Value = 683.9231
$ws.Range("J28").Value = 425
$ws.Range("K28").Value = 683.9231
$ws.Range("L28").Value = 425
$ws.Range("M28").Value = -198.9231
$ws.Range("N28").Value = -1395

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 47700
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 47700
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 47700
$ws.Range("N95").Value = -53192

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 966.5
$ws.Range("I135").Value = 759.8
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 6838.2
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -4303.2
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2558.625
$ws.Range("I138").Value = 1166.1666
$ws.Range("J138").Value = 3022.7778
$ws.Range("K138").Value = 3498.4998
$ws.Range("L138").Value = 9068.3334
$ws.Range("M138").Value = 1641.5002
$ws.Range("N138").Value = -19348.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 999
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -886
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 833
$ws.Range("I61").Value = 833
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 833
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 113331.664
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 113331.664
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 113331.664
$ws.Range("N92").Value = -118323.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 62000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 62000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 62000
$ws.Range("N95").Value = -67492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 999
$ws.Range("I116").Value = 999
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 999
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1295
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1552.75
$ws.Range("I132").Value = 1003.6667
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 3011.0001
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -481.0001000000002
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 833
$ws.Range("I136").Value = 833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2499
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 51

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 999
$ws.Range("I3").Value = 999
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -885
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2508.75
$ws.Range("I107").Value = 2510.5
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2510.5
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -590.5
$ws.Range("N107").Value = -6340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 36666.668
$ws.Range("I118").Value = 36666.668
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 36666.668
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -35009.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1316.3334
$ws.Range("I31").Value = 1316.3334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1316.3334
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1021.3334
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1316.3334
$ws.Range("I34").Value = 1316.3334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1316.3334
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1114.3334
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2436.625
$ws.Range("I58").Value = 2623.25
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 2623.25
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -2420.25
$ws.Range("N58").Value = -2656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1997.5
$ws.Range("I134").Value = 1997.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5992.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3457.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2436.625
$ws.Range("I136").Value = 2623.25
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 7869.75
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -5319.75
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 450.75
$ws.Range("I7").Value = 400.5
$ws.Range("J7").Value = 501
$ws.Range("K7").Value = 1201.5
$ws.Range("L7").Value = 1503
$ws.Range("M7").Value = -1089.5
$ws.Range("N7").Value = -1727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 400
$ws.Range("I122").Value = 400
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 450
$ws.Range("I80").Value = 450
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 450
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 548
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 450
$ws.Range("I83").Value = 450
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 2250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 2742
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1059.0769
$ws.Range("I107").Value = 766.125
$ws.Range("J107").Value = 1527.8
$ws.Range("K107").Value = 766.125
$ws.Range("L107").Value = 1527.8
$ws.Range("M107").Value = 1153.875
$ws.Range("N107").Value = -5367.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2922.6
$ws.Range("I132").Value = 2777
$ws.Range("J132").Value = 3505
$ws.Range("K132").Value = 8331
$ws.Range("L132").Value = 10515
$ws.Range("M132").Value = -5801
$ws.Range("N132").Value = -15575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 10000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 10000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 10225
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 10225
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 10225
$ws.Range("N104").Value = -17213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1202
$ws.Range("I132").Value = 1202
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3606
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1076

Write-Output "edit applied"
